# UF-6721 - Migrate PR3 data
# Adds a validation of legal id check digit ("kontrollsiffra") and fixes the
# stored PERSONNR values so the check digit is correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Correct the check digit ("kontrollsiffra") of each person's id number
# (column K / PERSONNR) on rows 2-4.
$ws.Range("K2").Value = "550717-7839"
$ws.Range("K3").Value = "420330-8947"
$ws.Range("K4").Value = "450627-8839"

# Leave the selection where the author left it when saving.
$null = $ws.Range("K11").Select()
